$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 151
$ws.Range("J2").Value = 730
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 181
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = 133
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 8
$ws.Range("S2").Value = 64
$ws.Range("T2").Value = 118
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 1100
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 1182
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 6
